$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.221.79'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '3.487.07'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '603.35'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').Value = '143.83'
$ws.Range('E6').Value = '  -2.83%  '
$ws.Range('D7').Value = '3.484.92'
$ws.Range('E7').Value = '  -0.40%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.476'
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range('D10').Value = '8.07'
$ws.Range('E10').Value = '  +0.40%  '
$ws.Range('E11').Value = '  -5.20%  '
$ws.Range('E12').Value = '  -2.86%  '
$ws.Range('D13').Value = '4.075.65'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('D14').Value = '30.33'
$ws.Range('E14').Value = '  -2.92%  '
$ws.Range('E15').Value = '  -5.48%  '
$ws.Range('D16').Value = '3.489.57'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('D17').Value = '66.233.82'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('D19').Value = '10.76'
$ws.Range('E19').Value = '  +3.82%  '
$ws.Range('E20').Value = '  -3.79%  '
$ws.Range('D21').Value = '14.80'
$ws.Range('E21').Value = '  -3.40%  '
$ws.Range('D22').Value = '425.61'
$ws.Range('E22').Value = '  -2.00%  '
$ws.Range('D23').Value = '0.595'
$ws.Range('E23').Value = '  -2.50%  '
$ws.Range('D24').Value = '77.84'
$ws.Range('E24').Value = '  -2.11%  '
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '3.622.34'
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').Value = '0.0000116'
$ws.Range('E27').Value = '  -2.90%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '9.24'
$ws.Range('E28').Value = '  -5.81%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '7.93'
$ws.Range('E29').Value = '  -3.77%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '2.46'
$ws.Range('E30').Value = '  -1.48%  '
$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').Value = '1.02'
$ws.Range('E31').Value = '  +2.12%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value = '0.165'
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value = '1.46'
$ws.Range('E33').Value = '  -8.83%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '25.08'
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('B35').Value = 'RenzoRestakedETH'
$ws.Range('C35').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D35').Value = '3.478.47'
$ws.Range('E35').Value = '  -0.38%  '
$ws.Range('B36').Value = 'USDe'
$ws.Range('C36').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '1.74'
$ws.Range('E37').Value = '  -3.63%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = '5.61'
$ws.Range('E38').Value = '  -5.60%  '
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').Value = '7.70'
$ws.Range('E39').Value = '  -3.55%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.10%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').Value = '169.80'
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').Value = '0.0858'
$ws.Range('E42').Value = '  -3.81%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = '5.15'
$ws.Range('E43').Value = '  -5.04%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = '0.879'
$ws.Range('E44').Value = '  -1.88%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').Value = '1.90'
$ws.Range('E45').Value = '  -9.20%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '45.38'
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '26.05'
$ws.Range('E47').Value = '  -10.69%  '
$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').Value = '1.20'
$ws.Range('E48').Value = '  -3.59%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').Value = '2.39'
$ws.Range('E49').Value = '  -1.04%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '7.13'
$ws.Range('E50').Value = '  -4.64%  '
$ws.Range('B51').Value = 'SuiNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D51').Value = '0.939'
$ws.Range('E51').Value = '  -3.06%  '
